$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2").Value = 46070
$ws.Range("C3").Value = 46070
$ws.Range("C4").Value = 46070
$ws.Range("C5").Value = 46070
$ws.Range("C6").Value = 46070

$ws.Range("A7").Value = "A 22256-2022"
$ws.Range("B7").Value = 44712
$ws.Range("C7").Value = 46070
$ws.Range("G7").Value = 11.2

$ws.Range("A8").Value = "A 46993-2025"
$ws.Range("B8").Value = 45929.54670138889
$ws.Range("C8").Value = 46070
$ws.Range("G8").Value = 2.8

$ws.Range("A9").Value = "A 46998-2025"
$ws.Range("B9").Value = 45929.54851851852
$ws.Range("C9").Value = 46070
$ws.Range("G9").Value = 0.9

$ws.Range("A10").Value = "A 57394-2024"
$ws.Range("B10").Value = 45629.6907175926
$ws.Range("C10").Value = 46070
$ws.Range("G10").Value = 0.5

$ws.Range("A11").Value = "A 63548-2025"
$ws.Range("B11").Value = 46013
$ws.Range("C11").Value = 46070

$ws.Range("A12").Value = "A 6679-2026"
$ws.Range("B12").Value = 46056.60961805555
$ws.Range("C12").Value = 46070
$ws.Range("G12").Value = 2.4

$ws.Range("A13").Value = "A 6684-2026"
$ws.Range("B13").Value = 46056.61989583333
$ws.Range("C13").Value = 46070
$ws.Range("G13").Value = 8.199999999999999

$ws.Range("A14").Value = "A 59877-2025"
$ws.Range("B14").Value = 45993
$ws.Range("C14").Value = 46070
$ws.Range("G14").Value = 1

$ws.Range("A15").Value = "A 59471-2024"
$ws.Range("B15").Value = 45638
$ws.Range("C15").Value = 46070

$ws.Range("A16").Value = "A 53218-2023"
$ws.Range("B16").Value = 45229
$ws.Range("C16").Value = 46070
$ws.Range("G16").Value = 5.4

$ws.Range("A17").Value = "A 20239-2025"
$ws.Range("B17").Value = 45772
$ws.Range("C17").Value = 46070
$ws.Range("G17").Value = 1.9

$ws.Range("A18").Value = "A 18118-2022"
$ws.Range("B18").Value = 44684
$ws.Range("C18").Value = 46070
$ws.Range("G18").Value = 4.2

$ws.Range("A19").Value = "A 57391-2024"
$ws.Range("B19").Value = 45629.68717592592
$ws.Range("C19").Value = 46070
$ws.Range("G19").Value = 1.8

$ws.Range("A20").Value = "A 54557-2023"
$ws.Range("B20").Value = 45233.6346875
$ws.Range("C20").Value = 46070
$ws.Range("G20").Value = 3.8

$ws.Range("A21").Value = "A 325-2024"
$ws.Range("B21").Value = 45295
$ws.Range("C21").Value = 46070
$ws.Range("G21").Value = 2.1

$ws.Range("A22").Value = "A 20054-2024"
$ws.Range("B22").Value = 45434
$ws.Range("C22").Value = 46070
$ws.Range("G22").Value = 7.3

$ws.Range("A23").Value = "A 26708-2023"
$ws.Range("B23").Value = 45093
$ws.Range("C23").Value = 46070
$ws.Range("G23").Value = 4.1

$ws.Range("C24").Value = 46070
$ws.Range("C25").Value = 46070
